$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" in H1, reusing the same formatting (style) as the
# other header cells (e.g. G1 "sum") by copying formats only.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("H1").Value = "Save"

# Fill H2:H15 with 0 (plain numeric cells, no special style)
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
